# outputLayout.xlsx update:
#   - insert a brand-new first worksheet "01_Change of Subjects" holding a
#     new "Question Code" list (12 new codes)
#   - keep the two existing worksheets ("20_Properties of Circles" and
#     "21_Prop of Tangent to Circle"), just shifted one slot to the right
#   - the new sheet becomes the active / selected tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet and fill it in BEFORE moving it - once it is
#    relocated to slot 1 the old object reference can go stale.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "01_Change of Subjects"

$questionCodes = @(
    "DSE24PII_Q03",
    "DSE21PII_Q02",
    "DSE20PII_Q02",
    "DSE13PII_Q02",
    "DSE16PII_Q02",
    "DSE22PII_Q05",
    "DSE18PII_Q02",
    "DSE23PII_Q01",
    "DSE19PII_Q05",
    "DSE17PII_Q03",
    "DSEPPPII_Q02",
    "DSESPPII_Q02"
)

$newSheet.Range("A1").Value = 1
$newSheet.Range("B1").Value = "Question Code"

for ($i = 0; $i -lt $questionCodes.Length; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $i + 1
    $newSheet.Cells.Item($r, 2).Value = $questionCodes[$i]
}

# Move the new sheet to the very front of the workbook.
$newSheet.Move($wb.Worksheets.Item(1))

# ---------------------------------------------------------------------
# 2. Re-resolve sheets by name (safest after the Move) and fix up the
#    active tab / selections to match the final layout.
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("01_Change of Subjects")
$sheet2 = $wb.Worksheets.Item("20_Properties of Circles")
$sheet3 = $wb.Worksheets.Item("21_Prop of Tangent to Circle")

# New first sheet becomes the active / visible tab with B10 selected.
$sheet1.Activate()
$sheet1.Range("B10").Select()

# Older sheets: refresh their selections (no longer the active tab).
$sheet2.Range("A1:B17").Select()
$sheet3.Range("A2:A14").Select()

# Land back on the new first sheet, matching the saved workbook state.
$sheet1.Activate()
